# Revert "change password feature is described":
# remove the two descriptive paragraphs that were added under the
# "Change password" bullet (the facility-description paragraph and the
# paragraph describing the submit/validate behaviour), restoring the
# bullet list so "Change password" is immediately followed by
# "User management".

$d = $word.ActiveDocument

$startMarker = "This facility allows the users to change their current password. In order to do this, they need to provide the current password and the new password."
$endMarker   = "When the values are submitted, system has to check if the old password of the user is valid or not. If valid, the new password will be updated in the database for the user. If the old password is not matching, a proper error message has to be displayed to the user. "

$paragraphs = $d.Paragraphs
$startIndex = -1
$endIndex = -1

for ($i = 1; $i -le $paragraphs.Count; $i++) {
    $p = $paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($startIndex -eq -1 -and $t -eq $startMarker) {
        $startIndex = $i
    }
    if ($endIndex -eq -1 -and $t -eq $endMarker) {
        $endIndex = $i
    }
}

if ($startIndex -ne -1 -and $endIndex -ne -1 -and $endIndex -ge $startIndex) {
    $startPara = $paragraphs.Item($startIndex)
    # Range up to (but excluding) the paragraph right after $endIndex so the
    # deletion also removes the paragraph marks of both target paragraphs.
    if ($endIndex -lt $paragraphs.Count) {
        $afterPara = $paragraphs.Item($endIndex + 1)
        $deleteRange = $d.Range($startPara.Range.Start, $afterPara.Range.Start)
    } else {
        $lastPara = $paragraphs.Item($endIndex)
        $deleteRange = $d.Range($startPara.Range.Start, $lastPara.Range.End)
    }
    $deleteRange.Delete()
}
